$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column headers are unchanged (row 1). Data rows 2-5 describe every
# sending/target cluster combination of FAPs and sCs for the Slurp1-Chrna7
# ligand/receptor pair.

# Row 2: FAPs -> FAPs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Slurp1"
$ws.Range("C2").Value = "Chrna7"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.4197276666666667
$ws.Range("H2").Value = 1.259183
$ws.Range("I2").Value = 0.6621316277815084
$ws.Range("J2").Value = 0.6621316277815084
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5188813333333333
$ws.Range("N2").Value = 1.556644
$ws.Range("O2").Value = 0.3951499920672696
$ws.Range("P2").Value = 0.3951499920672696
$ws.Range("Q2").Value = 0.2177888513168889
$ws.Range("R2").Value = 1.960099661852
$ws.Range("S2").Value = 0.2616413074653513
$ws.Range("T2").Value = 0.2616413074653513

# Row 3: FAPs -> sCs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Slurp1"
$ws.Range("C3").Value = "Chrna7"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.4197276666666667
$ws.Range("H3").Value = 1.259183
$ws.Range("I3").Value = 0.6621316277815084
$ws.Range("J3").Value = 0.6621316277815084
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7942436666666666
$ws.Range("N3").Value = 2.382731
$ws.Range("O3").Value = 0.6048500079327305
$ws.Range("P3").Value = 0.6048500079327305
$ws.Range("Q3").Value = 0.3333660409747777
$ws.Range("R3").Value = 3.000294368772999
$ws.Range("S3").Value = 0.4004903203161571
$ws.Range("T3").Value = 0.4004903203161571

# Row 4: sCs -> FAPs (new row)
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Slurp1"
$ws.Range("C4").Value = "Chrna7"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.214176
$ws.Range("H4").Value = 0.642528
$ws.Range("I4").Value = 0.3378683722184917
$ws.Range("J4").Value = 0.3378683722184917
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.5188813333333333
$ws.Range("N4").Value = 1.556644
$ws.Range("O4").Value = 0.3951499920672696
$ws.Range("P4").Value = 0.3951499920672696
$ws.Range("Q4").Value = 0.111131928448
$ws.Range("R4").Value = 1.000187356032
$ws.Range("S4").Value = 0.1335086846019183
$ws.Range("T4").Value = 0.1335086846019183

# Row 5: sCs -> sCs (new row)
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Slurp1"
$ws.Range("C5").Value = "Chrna7"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.214176
$ws.Range("H5").Value = 0.642528
$ws.Range("I5").Value = 0.3378683722184917
$ws.Range("J5").Value = 0.3378683722184917
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7942436666666666
$ws.Range("N5").Value = 2.382731
$ws.Range("O5").Value = 0.6048500079327305
$ws.Range("P5").Value = 0.6048500079327305
$ws.Range("Q5").Value = 0.170107931552
$ws.Range("R5").Value = 1.530971383968
$ws.Range("S5").Value = 0.2043596876165734
$ws.Range("T5").Value = 0.2043596876165734
